$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sending cluster changes from "FAPs" to "ECs"; target cluster stays "FAPs"
$ws.Range("A2").Value = "ECs"
$ws.Range("G2").Value = 0.03814
$ws.Range("H2").Value = 0.11442
$ws.Range("I2").Value = 0.0004360684493923871
$ws.Range("J2").Value = 0.0004360684493923871
$ws.Range("Q2").Value = 0.01047428649333333
$ws.Range("R2").Value = 0.09426857844
$ws.Range("S2").Value = 0.0004360684493923871
$ws.Range("T2").Value = 0.0004360684493923871

# Row 3: sending cluster changes from "MuSCs" to "FAPs"; target cluster stays "FAPs"
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2972863333333333
$ws.Range("H3").Value = 0.891859
$ws.Range("I3").Value = 0.003398982443686812
$ws.Range("J3").Value = 0.003398982443686811
$ws.Range("Q3").Value = 0.08164295295977778
$ws.Range("R3").Value = 0.734786576638
$ws.Range("S3").Value = 0.003398982443686812
$ws.Range("T3").Value = 0.003398982443686811

# Row 4: sending cluster stays "Resolving-Mac"; target cluster stays "FAPs"
$ws.Range("G4").Value = 87.12790666666666
$ws.Range("H4").Value = 261.38372
$ws.Range("I4").Value = 0.9961649491069209
$ws.Range("J4").Value = 0.9961649491069208
$ws.Range("Q4").Value = 23.92770466678222
$ws.Range("R4").Value = 215.34934200104
$ws.Range("S4").Value = 0.9961649491069209
$ws.Range("T4").Value = 0.9961649491069208
